$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_3_7_4", 0.3686433976866401, -0.01320909814009119, -4.101504988007466, -0.4466739788140937, 0.698725700378418, 0.8531020283699036, 0.5738874077796936, 0.7217071652412415),
    @("model_3_7_5", 0.3705493784890861, 0.03596296974446167, -4.416028201678258, -0.4361152690431001, 0.6966163516044617, 0.811700165271759, 0.609269380569458, 0.7164397239685059),
    @("model_3_7_6", 0.3751840325833462, 0.07336303098998997, -4.464094969199189, -0.4077961353898132, 0.691487193107605, 0.7802100777626038, 0.6146765351295471, 0.7023119926452637),
    @("model_3_7_9", 0.3758136480548316, 0.08002573490063103, -4.510481406227851, -0.4067645649212259, 0.6907902956008911, 0.7746002078056335, 0.6198947429656982, 0.701797366142273),
    @("model_3_7_10", 0.3780373620372752, 0.100909566770314, -4.563997677771314, -0.3937828850359306, 0.6883293390274048, 0.7570164203643799, 0.6259149312973022, 0.6953212022781372),
    @("model_3_7_7", 0.378107107187138, 0.09492858109353242, -4.47560572876631, -0.3897478260948024, 0.6882521510124207, 0.7620522975921631, 0.6159714460372925, 0.6933082342147827),
    @("model_3_7_11", 0.3789374662801361, 0.1088651615540336, -4.578295279176776, -0.3881914418912817, 0.6873332262039185, 0.7503179311752319, 0.6275233626365662, 0.6925317645072937),
    @("model_3_7_8", 0.3796500855543768, 0.1072547511145063, -4.473698810434736, -0.3785319648273051, 0.6865445375442505, 0.7516739368438721, 0.6157568693161011, 0.6877129077911377),
    @("model_3_7_3", 0.3799422210243737, 0.05535964063353349, -3.861543934984915, -0.3599421776360852, 0.6862212419509888, 0.7953685522079468, 0.546893298625946, 0.6784389019012451),
    @("model_3_7_13", 0.3818421508917313, 0.1363708721810055, -4.630844007983214, -0.3691929173361947, 0.6841185688972473, 0.7271586656570435, 0.6334347724914551, 0.6830539107322693),
    @("model_3_7_15", 0.3820027708219673, 0.1401857291422071, -4.659192134294292, -0.3687922909085832, 0.6839408278465271, 0.7239466309547424, 0.636623740196228, 0.6828539967536926),
    @("model_3_7_12", 0.3827403968779964, 0.1429537614240083, -4.615242574452959, -0.3616541774047737, 0.6831245422363281, 0.7216160297393799, 0.6316796541213989, 0.6792929768562317),
    @("model_3_7_17", 0.3830030732910701, 0.1533396747848375, -4.712695908752921, -0.3627157345105898, 0.6828337907791138, 0.712871253490448, 0.642642617225647, 0.6798226237297058),
    @("model_3_7_18", 0.3832732555695593, 0.1553890895183483, -4.712680673021319, -0.3608821245231937, 0.6825347542762756, 0.7111456990242004, 0.6426408290863037, 0.6789078712463379),
    @("model_3_7_19", 0.3837774447646137, 0.1593485079326379, -4.712320075153866, -0.3573068610524328, 0.6819767355918884, 0.7078119516372681, 0.6426002979278564, 0.6771242022514343),
    @("model_3_7_24", 0.3840521931831782, 0.1657461388625694, -4.754594924942715, -0.3560761554768501, 0.6816727519035339, 0.7024252414703369, 0.6473559141159058, 0.6765103340148926),
    @("model_3_7_20", 0.384352227256166, 0.1639538882617981, -4.712091648486959, -0.3531682164400634, 0.6813406348228455, 0.7039343118667603, 0.6425746083259583, 0.6750596165657043),
    @("model_3_7_23", 0.3852472039705374, 0.1749138050700778, -4.748098869888484, -0.3471945400800731, 0.6803502440452576, 0.6947062611579895, 0.6466251611709595, 0.6720795035362244),
    @("model_3_7_14", 0.3853622381308275, 0.1674220446012945, -4.651841163442049, -0.3436750008399758, 0.6802229285240173, 0.7010142207145691, 0.6357967853546143, 0.6703237295150757),
    @("model_3_7_22", 0.3853921067522368, 0.1753801031772197, -4.740251391987956, -0.3459460714044038, 0.6801897883415222, 0.694313645362854, 0.6457424163818359, 0.6714567542076111),
    @("model_3_7_21", 0.3856499151072567, 0.1769404972309508, -4.733684753118875, -0.3438539359620127, 0.6799045205116272, 0.6929998397827148, 0.6450036764144897, 0.6704129576683044),
    @("model_3_7_2", 0.3857241347267666, 0.09457723019048725, -3.727471631867465, -0.3106734320433244, 0.6798223853111267, 0.7623480558395386, 0.531810998916626, 0.6538600921630859),
    @("model_3_7_16", 0.3862736094839339, 0.1770572397193503, -4.6815391450447, -0.3382171602763813, 0.6792142987251282, 0.6929015517234802, 0.6391376256942749, 0.6676009297370911),
    @("model_3_7_0", 0.4229242061289731, 0.3516312149589934, -1.583865170870437, 0.1464801658792888, 0.6386528611183167, 0.5459136962890625, 0.2906686663627625, 0.4257983565330505),
    @("model_3_7_1", 0.4265960472444684, 0.3129803901681678, -2.099965243258279, 0.05717958372558618, 0.6345891952514648, 0.5784569978713989, 0.3487267196178436, 0.4703480005264282)

)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($j = 1; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}

Write-Host "Updated $($data.Count) rows"
